$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert a new row for the "shellyrelay" IO Command right before the
# existing "Hottop Heater" row (row 58), pushing all following rows down.
$ws.Rows.Item(58).Insert()

# Populate the new row: command usage in column B, description in column C
# (no entry in column A, matching the style of the preceding Kaleido rows).
$ws.Cells.Item(58, 2).Value = "shellyrelay(n,b)"
$ws.Cells.Item(58, 3).Value = "switches Shelly plug number <n> ON if b is true or 1, and OFF otherwise"

# Update the active selection to reflect the new location shown in the diff.
$ws.Range("C58").Select()
